$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price/volume figures (Price column is stored as text,
# so force a text NumberFormat before assigning the value and then
# clear the formatting override so the cell style matches the original).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.01"
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.00"
$ws.Range("D3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.203"
$ws.Range("D4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05783"
$ws.Range("D5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.508"
$ws.Range("D6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.121"
$ws.Range("D7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8149"
$ws.Range("D8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8490"
$ws.Range("D9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1358"
$ws.Range("D10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06963"
$ws.Range("D11").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02871"
$ws.Range("D13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09381"
$ws.Range("D14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.746"
$ws.Range("D15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001510"
$ws.Range("D16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04683"
$ws.Range("D17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005972"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006273"
$ws.Range("D19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001238"
$ws.Range("D20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004298"
$ws.Range("D21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008606"
$ws.Range("D22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.501"
$ws.Range("D23").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3172"
$ws.Range("D25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1338"
$ws.Range("D26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1326"
$ws.Range("D27").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03646"
$ws.Range("D40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006294"
$ws.Range("D41").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003102"
$ws.Range("D43").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005273"
$ws.Range("D45").ClearFormats()

$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002346"
$ws.Range("D48").ClearFormats()
